# Update the "Förändrad" (column C) date value for every data row
# from serial date 45204 (2023-10-05) to 45205 (2023-10-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (data starts at row 2).
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45205
    }
}
